$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# R1 reference value (F22) changed from 950000 to 98000
# ---------------------------------------------------------------------------
$ws.Range("F22").Value = 98000

# ---------------------------------------------------------------------------
# Drop the auxiliary "R1(Ohm)" label (I1) and its value (J1) -- the
# sharedStrings table will lose the now-unreferenced "R1(Ohm)" string on save.
# ---------------------------------------------------------------------------
$ws.Range("I1").ClearContents()
$ws.Range("J1").ClearContents()

# ---------------------------------------------------------------------------
# Column F (row 6..14): switch the explicit-base LOG() call for LOG10()
# ---------------------------------------------------------------------------
$ws.Range("F6").Formula = "=20*LOG10(B6*`$F`$22/C6)"
$ws.Range("F7:F14").Formula = "=20*LOG10(B7*`$F`$22/C7)"

# ---------------------------------------------------------------------------
# Column G (row 6..14): now holds the impedance ratio C/F22 instead of the
# raw B*F22/C product; also trim the stray shared formula in G15:G17 that
# used to spill #DIV/0! errors.
# ---------------------------------------------------------------------------
$ws.Range("G6").Formula = "=C6/`$F`$22"
$ws.Range("G7:G14").Formula = "=C7/`$F`$22"
$ws.Range("G15:G17").ClearContents()

# ---------------------------------------------------------------------------
# New column I (row 6..14): recompute |Z| dB from B and the new G ratio,
# matching F's number format (scientific, 2 decimals).
# ---------------------------------------------------------------------------
$ws.Range("I6").Formula = "=20*LOG(B6/G6)"
$ws.Range("I6").NumberFormat = "0.00E+00"
$ws.Range("I7:I14").Formula = "=20*LOG(B7/G7)"
$ws.Range("I7:I14").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Reposition/resize "Chart 1" (moved further right on the sheet).
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 711.0595703125
$co.Top = 27.75
$co.Width = 433.0625
$co.Height = 216

# ---------------------------------------------------------------------------
# Selection now sits on the recalculated impedance-ratio column.
# ---------------------------------------------------------------------------
$ws.Range("G6:G14").Select()
